# Add new columns to the approved_requests "Data" template:
#   - "Reseller ID #" inserted before "Reseller Name" (new col Q)
#   - "Vendor ID" and "Vendor Name" inserted after "HUB Name" (new cols Z:AA)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the "Reseller ID #" column at Q (shifts Reseller Name.. right by 1)
$ws.Columns("Q").Insert()
$ws.Range("Q1").Value = "Reseller ID #"

# 2) Insert two columns "Vendor ID" / "Vendor Name" at Z:AA (after the now-shifted HUB Name column Y)
$ws.Columns("Z:AA").Insert()
$ws.Range("Z1").Value = "Vendor ID"
$ws.Range("AA1").Value = "Vendor Name"

# 3) Re-apply the AutoFilter over the new full header range
$ws.AutoFilterMode = $false
$ws.Range("A1:AK1").AutoFilter()

# 4) Update the hidden _FilterDatabase defined name to the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AK`$1"
    }
}

# 5) Match the widened "Reseller External ID" column width
$ws.Columns("P").ColumnWidth = 21.36328125
$ws.Columns("Q").ColumnWidth = 18.6328125

# 6) Restore the selected cell shown in the workbook
$ws.Range("AA2").Select()
